$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rotation_req")

$ws.Cells.Item(1, 1).Value = "AAAAAa"
$ws.Cells.Item(1, 2).Value = "AAAAA"
$ws.Cells.Item(1, 3).Value = 1
$ws.Cells.Item(2, 1).Value = "AAAAAa"
$ws.Cells.Item(2, 2).Value = "GGGGA"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 1).Value = "GGGGAa2"
$ws.Cells.Item(3, 2).Value = "GGGGA"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(4, 1).Value = "GGGANa2"
$ws.Cells.Item(4, 2).Value = "GGGAN"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(5, 1).Value = "GGGAEa2"
$ws.Cells.Item(5, 2).Value = "GGGAE"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 1).Value = "GGGAOFa2"
$ws.Cells.Item(6, 2).Value = "GGGAE"
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(7, 1).Value = "GGGAOFa2"
$ws.Cells.Item(7, 2).Value = "GGGAOF"
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(8, 1).Value = "GGGC1Na2"
$ws.Cells.Item(8, 2).Value = "GGGC1N"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(9, 1).Value = "GGGC1Ea2"
$ws.Cells.Item(9, 2).Value = "GGGC1E"
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(10, 1).Value = "GGGC1Pa2"
$ws.Cells.Item(10, 2).Value = "GGGC1P"
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(11, 1).Value = "GGGC1OFa2"
$ws.Cells.Item(11, 2).Value = "GGGC1E"
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 1).Value = "GGGC1OFa2"
$ws.Cells.Item(12, 2).Value = "GGGC1OF"
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(13, 1).Value = "SPSPSPSPSPsp"
$ws.Cells.Item(13, 2).Value = "SPSPSPSPSP"
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(14, 1).Value = "YABABa"
$ws.Cells.Item(14, 2).Value = "GGGAE"
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(15, 1).Value = "YABABa"
$ws.Cells.Item(15, 2).Value = "YABAB"
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(16, 1).Value = "AYABAb"
$ws.Cells.Item(16, 2).Value = "AYABA"
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(17, 1).Value = "AYABAb"
$ws.Cells.Item(17, 2).Value = "GGGGA"
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(18, 1).Value = "AYABAbd"
$ws.Cells.Item(18, 2).Value = "AYABA"
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(19, 1).Value = "AYABAbd"
$ws.Cells.Item(19, 2).Value = "GGGGA"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(20, 1).Value = "YAOAOFa"
$ws.Cells.Item(20, 2).Value = "GGGAE"
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(21, 1).Value = "YAOAOFa"
$ws.Cells.Item(21, 2).Value = "GGGAOF"
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(22, 1).Value = "YAOAOFa"
$ws.Cells.Item(22, 2).Value = "YAOAOF"
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(23, 1).Value = "AYAOAof"
$ws.Cells.Item(23, 2).Value = "AYAOA"
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(24, 1).Value = "AYAOAof"
$ws.Cells.Item(24, 2).Value = "GGGGA"
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(25, 1).Value = "YYNBNb"
$ws.Cells.Item(25, 2).Value = "GGGC1N"
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(26, 1).Value = "YYNBNb"
$ws.Cells.Item(26, 2).Value = "YYNBN"
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(27, 1).Value = "YYBNBz"
$ws.Cells.Item(27, 2).Value = "GGGC1E"
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(28, 1).Value = "YYBNBz"
$ws.Cells.Item(28, 2).Value = "YYBNB"
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(29, 1).Value = "YYBNBr"
$ws.Cells.Item(29, 2).Value = "GGGC1E"
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(30, 1).Value = "YYBNBr"
$ws.Cells.Item(30, 2).Value = "YYBNB"
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(31, 1).Value = "YYNWBr"
$ws.Cells.Item(31, 2).Value = "GGGC1E"
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(32, 1).Value = "YYNWBr"
$ws.Cells.Item(32, 2).Value = "YYNWB"
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(33, 1).Value = "YYWBNw"
$ws.Cells.Item(33, 2).Value = "GGGC1N"
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(34, 1).Value = "YYWBNw"
$ws.Cells.Item(34, 2).Value = "YYWBN"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(35, 1).Value = "YYBNWb"
$ws.Cells.Item(35, 2).Value = "GGGC1E"
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(36, 1).Value = "YYBNWb"
$ws.Cells.Item(36, 2).Value = "YYBNW"
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(37, 1).Value = "AAAAAz"
$ws.Cells.Item(37, 2).Value = "AAAAA"
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(38, 1).Value = "AAAAAz"
$ws.Cells.Item(38, 2).Value = "GGGGA"
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(39, 1).Value = "AAAANw"
$ws.Cells.Item(39, 2).Value = "AAAAN"
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(40, 1).Value = "AAAANw"
$ws.Cells.Item(40, 2).Value = "GGGAN"
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(41, 1).Value = "AAANWb"
$ws.Cells.Item(41, 2).Value = "AAANW"
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(42, 1).Value = "AAANWb"
$ws.Cells.Item(42, 2).Value = "GGGC1E"
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(43, 1).Value = "AANWBa"
$ws.Cells.Item(43, 2).Value = "AANWB"
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(44, 1).Value = "AANWBa"
$ws.Cells.Item(44, 2).Value = "GGGC1E"
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(45, 1).Value = "AYWBAa"
$ws.Cells.Item(45, 2).Value = "AYWBA"
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(46, 1).Value = "AYWBAa"
$ws.Cells.Item(46, 2).Value = "GGGGA"
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(47, 1).Value = "YYBAAa"
$ws.Cells.Item(47, 2).Value = "GGGGA"
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(48, 1).Value = "YYBAAa"
$ws.Cells.Item(48, 2).Value = "YYBAA"
$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(49, 1).Value = "YYAAAa"
$ws.Cells.Item(49, 2).Value = "GGGGA"
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(50, 1).Value = "YYAAAa"
$ws.Cells.Item(50, 2).Value = "YYAAA"
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(51, 1).Value = "YAAAAa"
$ws.Cells.Item(51, 2).Value = "GGGGA"
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(52, 1).Value = "YAAAAa"
$ws.Cells.Item(52, 2).Value = "YAAAA"
$ws.Cells.Item(52, 3).Value = 1
$ws.Cells.Item(53, 1).Value = "YYFWBz"
$ws.Cells.Item(53, 2).Value = "GGGC1E"
$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(54, 1).Value = "YYFWBz"
$ws.Cells.Item(54, 2).Value = "YYFWB"
$ws.Cells.Item(54, 3).Value = 1
$ws.Cells.Item(55, 1).Value = "YYBNWf"
$ws.Cells.Item(55, 2).Value = "GGGC1E"
$ws.Cells.Item(55, 3).Value = 1
$ws.Cells.Item(56, 1).Value = "YYBNWf"
$ws.Cells.Item(56, 2).Value = "YYBNW"
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(57, 1).Value = "YYNWFw"
$ws.Cells.Item(57, 2).Value = "GGGC1P"
$ws.Cells.Item(57, 3).Value = 1
$ws.Cells.Item(58, 1).Value = "YYNWFw"
$ws.Cells.Item(58, 2).Value = "YYNWF"
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(59, 1).Value = "YYWFWb"
$ws.Cells.Item(59, 2).Value = "GGGC1E"
$ws.Cells.Item(59, 3).Value = 1
$ws.Cells.Item(60, 1).Value = "YYWFWb"
$ws.Cells.Item(60, 2).Value = "YYWFW"
$ws.Cells.Item(60, 3).Value = 1
$ws.Cells.Item(61, 1).Value = "YYWNWl"
$ws.Cells.Item(61, 2).Value = "GGGC1E"
$ws.Cells.Item(61, 3).Value = 1
$ws.Cells.Item(62, 1).Value = "YYWNWl"
$ws.Cells.Item(62, 2).Value = "YYWNW"
$ws.Cells.Item(62, 3).Value = 1
$ws.Cells.Item(63, 1).Value = "YYLWNw"
$ws.Cells.Item(63, 2).Value = "GGGC1N"
$ws.Cells.Item(63, 3).Value = 1
$ws.Cells.Item(64, 1).Value = "YYLWNw"
$ws.Cells.Item(64, 2).Value = "YYLWN"
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(65, 1).Value = "YYWLWz"
$ws.Cells.Item(65, 2).Value = "GGGC1E"
$ws.Cells.Item(65, 3).Value = 1
$ws.Cells.Item(66, 1).Value = "YYWLWz"
$ws.Cells.Item(66, 2).Value = "YYWLW"
$ws.Cells.Item(66, 3).Value = 1
$ws.Cells.Item(67, 1).Value = "YYNWLw"
$ws.Cells.Item(67, 2).Value = "GGGC1P"
$ws.Cells.Item(67, 3).Value = 1
$ws.Cells.Item(68, 1).Value = "YYNWLw"
$ws.Cells.Item(68, 2).Value = "YYNWL"
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(69, 1).Value = "AYAAAw"
$ws.Cells.Item(69, 2).Value = "AYAAA"
$ws.Cells.Item(69, 3).Value = 1
$ws.Cells.Item(70, 1).Value = "AYAAAw"
$ws.Cells.Item(70, 2).Value = "GGGGA"
$ws.Cells.Item(70, 3).Value = 1
$ws.Cells.Item(71, 1).Value = "AAWAAa"
$ws.Cells.Item(71, 2).Value = "AAWAA"
$ws.Cells.Item(71, 3).Value = 1
$ws.Cells.Item(72, 1).Value = "AAWAAa"
$ws.Cells.Item(72, 2).Value = "GGGGA"
$ws.Cells.Item(72, 3).Value = 1
$ws.Cells.Item(73, 1).Value = "AAAWAa"
$ws.Cells.Item(73, 2).Value = "AAAWA"
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(74, 1).Value = "AAAWAa"
$ws.Cells.Item(74, 2).Value = "GGGGA"
$ws.Cells.Item(74, 3).Value = 1
$ws.Cells.Item(75, 1).Value = "YAAAWa"
$ws.Cells.Item(75, 2).Value = "GGGAE"
$ws.Cells.Item(75, 3).Value = 1
$ws.Cells.Item(76, 1).Value = "YAAAWa"
$ws.Cells.Item(76, 2).Value = "YAAAW"
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(77, 1).Value = "YAAANw"
$ws.Cells.Item(77, 2).Value = "GGGAN"
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(78, 1).Value = "YAAANw"
$ws.Cells.Item(78, 2).Value = "YAAAN"
$ws.Cells.Item(78, 3).Value = 1
$ws.Cells.Item(79, 1).Value = "AYAAAz"
$ws.Cells.Item(79, 2).Value = "AYAAA"
$ws.Cells.Item(79, 3).Value = 1
$ws.Cells.Item(80, 1).Value = "AYAAAz"
$ws.Cells.Item(80, 2).Value = "GGGGA"
$ws.Cells.Item(80, 3).Value = 1
$ws.Cells.Item(81, 1).Value = "AYWAAa"
$ws.Cells.Item(81, 2).Value = "AYWAA"
$ws.Cells.Item(81, 3).Value = 1
$ws.Cells.Item(82, 1).Value = "AYWAAa"
$ws.Cells.Item(82, 2).Value = "GGGGA"
$ws.Cells.Item(82, 3).Value = 1
$ws.Cells.Item(83, 1).Value = "AANWAa"
$ws.Cells.Item(83, 2).Value = "AANWA"
$ws.Cells.Item(83, 3).Value = 1
$ws.Cells.Item(84, 1).Value = "AANWAa"
$ws.Cells.Item(84, 2).Value = "GGGGA"
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(85, 1).Value = "AAANWa"
$ws.Cells.Item(85, 2).Value = "AAANW"
$ws.Cells.Item(85, 3).Value = 1
$ws.Cells.Item(86, 1).Value = "AAANWa"
$ws.Cells.Item(86, 2).Value = "GGGC1E"
$ws.Cells.Item(86, 3).Value = 1
